$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "Datos actualizados" timestamp (row 1) ---
$ws.Range("A1").Value = "Datos actualizados a 27 de Marzo de 2020 a las 13:42"

# --- Move "Caceres" up: new figures, now placed right after "Gran Canaria" ---
# Insert a new row right after Gran Canaria (row 20) and populate it with the
# updated Caceres data.
$ws.Rows(21).Insert()
$ws.Range("A21").Value = "Caceres"
$ws.Range("B21").Value = 841
$ws.Range("C21").Value = 4
$ws.Range("D21").Value = 777
$ws.Range("E21").Value = 60

# Remove the old Caceres row, which has now shifted down to row 28
# (A28 = Caceres, 640, 4, 587, 35).
$ws.Rows(28).Delete()

# --- Move "Badajoz" up: new figures, now placed right after "Castello/Castellon" ---
# After the edits above, Castello/Castellon sits at row 35, so insert the new
# row right after it.
$ws.Rows(36).Insert()
$ws.Range("A36").Value = "Badajoz"
$ws.Range("B36").Value = 390
$ws.Range("C36").Value = 29
$ws.Range("D36").Value = 352
$ws.Range("E36").Value = 9

# Remove the old Badajoz row, which has now shifted down to row 41
# (A41 = Badajoz, 329, 20, 316, 4).
$ws.Rows(41).Delete()
